$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lookup table of rows whose Notified/Actual Production values changed.
# Key = row number (string), Value = @(newB, newC) where $null means "unchanged".
$bcChanges = @{
    "6" = @(0.29, $null)
    "14" = @(0.45, $null)
    "18" = @(0, $null)
    "22" = @(0.518, $null)
    "23" = @(0.51, $null)
    "24" = @(0.542, $null)
    "25" = @(0.631, $null)
    "26" = @(1.574, $null)
    "27" = @(1.741, $null)
    "28" = @(2.591, $null)
    "29" = @(6.69, $null)
    "30" = @(60.031, 12)
    "31" = @(108.673, 42)
    "32" = @(182.639, 89)
    "33" = @(262.772, 150)
    "34" = @(483.322, 211)
    "35" = @(592.737, 274)
    "36" = @(705.421, 344)
    "37" = @(825.7190000000001, 418)
    "38" = @(1070.307, 472)
    "39" = @(1175.008, 534)
    "40" = @(1284.423, 587)
    "41" = @(1363.79, 636)
    "42" = @(1517.691, 678)
    "43" = @(1581.554, 714)
    "44" = @(1634.925, 755)
    "45" = @(1675.83, 817)
    "46" = @(1719.85, 839)
    "47" = @(1732.729, 846)
    "48" = @(1739.835, 852)
    "49" = @(1730.922, 876)
    "50" = @(1710.712, 847)
    "51" = @(1680.604, 831)
    "52" = @(1636.304, 803)
    "53" = @(1581.226, 778)
    "54" = @(1451.575, $null)
    "55" = @(1363.08, 713)
    "56" = @(1260.199, 631)
    "57" = @(1137.318, 557)
    "58" = @(892.875, 498)
    "59" = @(774.164, 439)
    "60" = @(645.164, 350)
    "61" = @(524.711, 269)
    "62" = @(304.876, 196)
    "63" = @(216.257, 132)
    "64" = @(131.209, 73)
    "65" = @(77.193, 24)
    "66" = @(22.572, $null)
    "67" = @(15.088, $null)
    "68" = @(10.812, $null)
    "69" = @(10.554, $null)
    "70" = @(3.016, $null)
    "73" = @(2.676, $null)
    "74" = @(2.65, $null)
    "75" = @(0.65, $null)
    "76" = @(0, $null)
    "78" = @(0.49, $null)
    "82" = @(0.65, $null)
    "86" = @(0.45, $null)
    "102" = @(0.29, $null)
    "110" = @(0.45, $null)
    "114" = @(0, $null)
    "118" = @(0.462, $null)
    "119" = @(0.47, $null)
    "120" = @(0, $null)
    "121" = @(0.484, $null)
    "122" = @(0.588, $null)
    "123" = @(0.841, $null)
    "124" = @(1.197, $null)
    "125" = @(2.256, $null)
    "126" = @(12.012, 1)
    "127" = @(24.118, 0)
    "128" = @(39.222, 0)
    "129" = @(56.534, 0)
    "130" = @(93.496, 0)
    "131" = @(116.008, 0)
    "132" = @(142.186, 0)
    "133" = @(166.037, 0)
    "134" = @(206.921, 0)
    "135" = @(227.295, 0)
    "136" = @(249.076, $null)
    "137" = @(265.005, $null)
    "138" = @(292.77, $null)
    "139" = @(307.295, $null)
    "140" = @(327.611, $null)
    "141" = @(337.767, $null)
    "142" = @(353.525, $null)
    "143" = @(357.384, $null)
    "144" = @(357.203, $null)
    "145" = @(356.396, $null)
    "146" = @(353.526, $null)
    "147" = @(344.782, $null)
    "148" = @(334.905, $null)
    "149" = @(320.995, $null)
    "150" = @(301.379, $null)
    "151" = @(282.196, $null)
    "152" = @(260.077, $null)
    "153" = @(237.102, $null)
    "154" = @(191.808, $null)
    "155" = @(165.562, $null)
    "156" = @(136.962, $null)
    "157" = @(112.947, $null)
    "158" = @(76.873, $null)
    "159" = @(56.2, $null)
    "160" = @(39.025, $null)
    "161" = @(26.44, $null)
    "162" = @(13.986, $null)
    "163" = @(10.226, $null)
    "164" = @(8.977, $null)
    "165" = @(8.711, $null)
    "166" = @(2.65, $null)
    "169" = @(0, $null)
    "170" = @(0, $null)
    "171" = @(0.65, $null)
    "172" = @(0, $null)
    "174" = @(0.49, $null)
    "178" = @(0.65, $null)
    "182" = @(0.45, $null)
}

$dayShift = 22

for ($r = 2; $r -le 193; $r++) {
    # Column A: Timestamp -> shift the underlying date serial by +22 days
    $oldSerial = $ws.Cells.Item($r, 1).Value2
    $newSerial = $oldSerial + $dayShift
    $ws.Cells.Item($r, 1).Value2 = $newSerial

    # Column D: Quarter index stays the same, used to rebuild the Lookup text
    $quarter = $ws.Cells.Item($r, 4).Value2

    # Column E: Lookup = "DD.MM.YYYY" & Quarter, derived from the new date
    $newDate = [DateTime]::FromOADate($newSerial)
    $dayStr = $newDate.Day.ToString().PadLeft(2,'0')
    $monStr = $newDate.Month.ToString().PadLeft(2,'0')
    $yearStr = $newDate.Year.ToString()
    $ws.Cells.Item($r, 5).Value2 = "$dayStr.$monStr.$yearStr" + [string]$quarter

    # Columns B/C: apply any explicit value change for this row
    $key = [string]$r
    if ($bcChanges.ContainsKey($key)) {
        $pair = $bcChanges[$key]
        if ($null -ne $pair[0]) { $ws.Cells.Item($r, 2).Value2 = $pair[0] }
        if ($null -ne $pair[1]) { $ws.Cells.Item($r, 3).Value2 = $pair[1] }
    }
}

Write-Output "done"
